$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Learning Module column (C) cells 9-12 to include markdown links
$ws.Range("C9").Value = "[Information Processing](https://www.crumplab.com/cognition/articles/modules/L6_Information_processing.html)"
$ws.Range("C10").Value = "[Memory I](https://www.crumplab.com/cognition/articles/modules/L7_Memory_I.html)"
$ws.Range("C11").Value = "[Memory II](https://www.crumplab.com/cognition/articles/modules/L8_Memory_II.html)"
$ws.Range("C12").Value = "[Implicit Cognition](https://www.crumplab.com/cognition/articles/modules/L9_Implicit_Cognition.html) "

# Row 12 height increases to accommodate the longer wrapped text
$ws.Rows.Item(12).RowHeight = 68

# Update the view selection/scroll position to reflect where the edit took place
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C12").Select()
